# Player Performance workbook update:
#  1. Insert a new "Player Info" worksheet before the existing "ODI Batting"
#     sheet, with player bio data.
#  2. On the "ODI Batting" sheet, rename the MATCH_CARD_LINK column to
#     MATCH_CODE and replace the full scorecard URLs with just the numeric
#     match code.

$wb = $excel.ActiveWorkbook

# --- Capture a handle to the existing ("ODI Batting") sheet before adding
#     the new one, then re-resolve it by name afterwards -- inserting a
#     sheet ahead of it shifts its position, and we want a fresh handle.
$battingName = $wb.Worksheets.Item(1).Name

# --- Step 1: create the new "Player Info" sheet positioned before ODI Batting
$infoWs = $wb.Worksheets.Add($wb.Worksheets.Item(1))
$infoWs.Name = "Player Info"

$battingWs = $wb.Worksheets.Item($battingName)

# Copy the header row's cell formatting (bold, centered, bordered) from the
# ODI Batting header onto the new sheet's header row, then overwrite the
# text -- this reuses the existing header style instead of creating a new one.
$battingWs.Range("A1:D1").Copy($infoWs.Range("A1:D1"))

$infoWs.Range("A1").Value = "ID"
$infoWs.Range("B1").Value = "NAME"
$infoWs.Range("C1").Value = "BATTING_HAND"
$infoWs.Range("D1").Value = "BOWL_STYLE"

$infoWs.Range("A2").NumberFormat = "@"
$infoWs.Range("A2").Value = "4759"
$infoWs.Range("B2").Value = "Prithvi Pankaj Shaw"
$infoWs.Range("C2").Value = "Right Handed"
$infoWs.Range("D2").Value = "Right Arm Off Break"

# --- Step 2: update ODI Batting sheet - MATCH_CARD_LINK -> MATCH_CODE
$battingWs.Range("D1").Value = "MATCH_CODE"

$battingWs.Range("D2").NumberFormat = "@"
$battingWs.Range("D2").Value = "4402"

$battingWs.Range("D3").NumberFormat = "@"
$battingWs.Range("D3").Value = "4406"

$battingWs.Range("D4").NumberFormat = "@"
$battingWs.Range("D4").Value = "4410"

$battingWs.Range("D5").NumberFormat = "@"
$battingWs.Range("D5").Value = "4480"

$battingWs.Range("D6").NumberFormat = "@"
$battingWs.Range("D6").Value = "4482"

$battingWs.Range("D7").NumberFormat = "@"
$battingWs.Range("D7").Value = "4485"
